# Auto-generated edit script: updates Excalibur_Profits market-board derived
# columns H-N across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 75477.37
$ws.Range("J87").Value = 90781.375
$ws.Range("L87").Value = 90781.375
$ws.Range("N87").Value = -93277.375
$ws.Range("H90").Value = 75477.37
$ws.Range("J90").Value = 90781.375
$ws.Range("L90").Value = 272344.125
$ws.Range("N90").Value = -284824.125
$ws.Range("H109").Value = 70000
$ws.Range("I109").Value = 70000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 70000
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -68613
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 1530.1951
$ws.Range("I112").Value = 845
$ws.Range("J112").Value = 1584.2894
$ws.Range("K112").Value = 2535
$ws.Range("L112").Value = 4752.8682
$ws.Range("M112").Value = -1427
$ws.Range("N112").Value = -6968.8682
$ws.Range("H113").Value = 4068.5
$ws.Range("I113").Value = 4725
$ws.Range("J113").Value = 3630.8333
$ws.Range("K113").Value = 4725
$ws.Range("L113").Value = 3630.8333
$ws.Range("M113").Value = -1471
$ws.Range("N113").Value = -10138.8333
$ws.Range("H116").Value = 39004.65
$ws.Range("I116").Value = 46421.465
$ws.Range("J116").Value = 16754.2
$ws.Range("K116").Value = 46421.465
$ws.Range("L116").Value = 16754.2
$ws.Range("M116").Value = -42979.465
$ws.Range("N116").Value = -23638.2
$ws.Range("H127").Value = 3905.3333
$ws.Range("I127").Value = 4596.4
$ws.Range("K127").Value = 13789.2
$ws.Range("M127").Value = -8829.199999999999
$ws.Range("H138").Value = 2699.102
$ws.Range("I138").Value = 1794.84
$ws.Range("J138").Value = 3641.0417
$ws.Range("K138").Value = 5384.52
$ws.Range("L138").Value = 10923.1251
$ws.Range("M138").Value = -244.5199999999995
$ws.Range("N138").Value = -21203.1251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5249.452
$ws.Range("I32").Value = 4672.919
$ws.Range("K32").Value = 4672.919
$ws.Range("M32").Value = -4385.919
$ws.Range("H74").Value = 3879.8125
$ws.Range("I74").Value = 3620.7046
$ws.Range("K74").Value = 3620.7046
$ws.Range("M74").Value = -2746.7046
$ws.Range("H77").Value = 3879.8125
$ws.Range("I77").Value = 3620.7046
$ws.Range("K77").Value = 18103.523
$ws.Range("M77").Value = -13735.523
$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960
$ws.Range("H139").Value = 84999.836
$ws.Range("J139").Value = 84999.836
$ws.Range("L139").Value = 84999.836
$ws.Range("N139").Value = -95279.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 174127.83
$ws.Range("J140").Value = 174127.83
$ws.Range("L140").Value = 174127.83
$ws.Range("N140").Value = -184487.83

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 893.8889
$ws.Range("I16").Value = 578
$ws.Range("J16").Value = 1999.5
$ws.Range("K16").Value = 578
$ws.Range("L16").Value = 1999.5
$ws.Range("M16").Value = -291
$ws.Range("N16").Value = -2573.5
$ws.Range("H31").Value = 5207
$ws.Range("I31").Value = 2532.4375
$ws.Range("K31").Value = 2532.4375
$ws.Range("M31").Value = -2237.4375
$ws.Range("H34").Value = 5207
$ws.Range("I34").Value = 2532.4375
$ws.Range("K34").Value = 2532.4375
$ws.Range("M34").Value = -2330.4375
$ws.Range("H113").Value = 893.8889
$ws.Range("I113").Value = 578
$ws.Range("J113").Value = 1999.5
$ws.Range("K113").Value = 578
$ws.Range("L113").Value = 1999.5
$ws.Range("M113").Value = 1592
$ws.Range("N113").Value = -6339.5
$ws.Range("H122").Value = 3996.3
$ws.Range("I122").Value = 1525
$ws.Range("K122").Value = 4575
$ws.Range("M122").Value = -2125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3427.1428
$ws.Range("I5").Value = 665
$ws.Range("J5").Value = 20000
$ws.Range("K5").Value = 1995
$ws.Range("L5").Value = 60000
$ws.Range("M5").Value = -1883
$ws.Range("N5").Value = -60224
$ws.Range("H57").Value = 4403
$ws.Range("J57").Value = 6667.143
$ws.Range("L57").Value = 20001.429
$ws.Range("N57").Value = -21119.429
$ws.Range("H107").Value = 791.9
$ws.Range("J107").Value = 856.1667
$ws.Range("L107").Value = 2568.5001
$ws.Range("N107").Value = -6408.5001
$ws.Range("H118").Value = 899.3333
$ws.Range("I118").Value = 799
$ws.Range("K118").Value = 2397
$ws.Range("M118").Value = -1154
$ws.Range("H131").Value = 1689.3
$ws.Range("J131").Value = 1941
$ws.Range("L131").Value = 5823
$ws.Range("N131").Value = -15903
$ws.Range("H135").Value = 3427.1428
$ws.Range("I135").Value = 665
$ws.Range("J135").Value = 20000
$ws.Range("K135").Value = 5985
$ws.Range("L135").Value = 180000
$ws.Range("M135").Value = -3450
$ws.Range("N135").Value = -185070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27996.6
$ws.Range("I70").Value = 54994.5
$ws.Range("K70").Value = 54994.5
$ws.Range("M70").Value = -54724.5
$ws.Range("H73").Value = 27996.6
$ws.Range("I73").Value = 54994.5
$ws.Range("K73").Value = 54994.5
$ws.Range("M73").Value = -54058.5
$ws.Range("H102").Value = 5030.5312
$ws.Range("I102").Value = 5195.731
$ws.Range("K102").Value = 5195.731
$ws.Range("M102").Value = -3573.731
$ws.Range("H126").Value = 3036.8276
$ws.Range("I126").Value = 3002.5715
$ws.Range("J126").Value = 3047.7273
$ws.Range("K126").Value = 9007.7145
$ws.Range("L126").Value = 9143.1819
$ws.Range("M126").Value = -6537.7145
$ws.Range("N126").Value = -14083.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3040.6
$ws.Range("I7").Value = 2232.6667
$ws.Range("K7").Value = 2232.6667
$ws.Range("M7").Value = -2120.6667
$ws.Range("H93").Value = 2355.818
$ws.Range("I93").Value = 2008.5
$ws.Range("K93").Value = 2008.5
$ws.Range("M93").Value = -760.5
$ws.Range("H126").Value = 3040.6
$ws.Range("I126").Value = 2232.6667
$ws.Range("K126").Value = 6698.000100000001
$ws.Range("M126").Value = -4228.000100000001
$ws.Range("H132").Value = 3620
$ws.Range("I132").Value = 3277.4
$ws.Range("K132").Value = 9832.200000000001
$ws.Range("M132").Value = -7302.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 78742.664
$ws.Range("J75").Value = 78742.664
$ws.Range("L75").Value = 78742.664
$ws.Range("N75").Value = -80614.664
$ws.Range("H78").Value = 78742.664
$ws.Range("J78").Value = 78742.664
$ws.Range("L78").Value = 236227.992
$ws.Range("N78").Value = -245587.992
$ws.Range("H126").Value = 5538.5557
$ws.Range("I126").Value = 5606
$ws.Range("K126").Value = 16818
$ws.Range("M126").Value = -14348
$ws.Range("H136").Value = 8641739
$ws.Range("I136").Value = 12800625
$ws.Range("K136").Value = 38401875
$ws.Range("M136").Value = -38399325
